$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python")

# Row 33 gains a Level (column D) value that was previously blank
$ws.Range("D33").Value = "Medium"

# New row 34: Different Ways to Add Parentheses
$ws.Range("A34").Value = "？"
$ws.Range("B34").Value = "Different Ways to Add Parentheses"
$ws.Range("C34").Value = 241
$ws.Range("D34").Value = "Medium"

# New row 35: Prime Arrangements
$ws.Range("A35").Value = "Y"
$ws.Range("B35").Value = " Prime Arrangements"
$ws.Range("C35").Value = 1175
$ws.Range("D35").Value = "Easy"

# New row 36: Valid Parentheses
$ws.Range("A36").Value = "Y"
$ws.Range("B36").Value = "Valid Parentheses"
$ws.Range("C36").Value = 20
$ws.Range("D36").Value = "Easy"

# Match the new active selection recorded in the workbook
[void]$ws.Range("G36").Select()
